# Fruta / hortaliza, semanal
# New weekly reading arrives; insert it as the new first data row (row 2)
# and push the existing rows down by one (row 2 -> row 3, ... row 13 -> row 14).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 2 (the existing data moves
# down to make room for the new weekly observation).
$ws.Rows.Item(2).Insert()

# The insert copies the formatting of the row above (the bold header row),
# so strip that back to a plain, unstyled data row like all the others.
$ws.Range("A2:T2").ClearFormats()

# Column D carries the date/time number format used by every other row in
# this column - copy it from the row right below (the row that used to be
# row 2) instead of hard-coding a style index.
$ws.Cells.Item(2, 4).NumberFormat = $ws.Cells.Item(3, 4).NumberFormat

# Fill in the constant columns (same market/product/category on every row).
$ws.Cells.Item(2, 1).Value = 7
$ws.Cells.Item(2, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(2, 3).Value = "Ñuble"
$ws.Cells.Item(2, 4).Value = 45043
$ws.Cells.Item(2, 5).Value = 16
$ws.Cells.Item(2, 6).Value = "Fruta"
$ws.Cells.Item(2, 7).Value = 100107
$ws.Cells.Item(2, 8).Value = "Otros"
$ws.Cells.Item(2, 9).Value = 100107011
$ws.Cells.Item(2, 10).Value = "Tuna"
$ws.Cells.Item(2, 11).Value = "Sin especificar"
$ws.Cells.Item(2, 12).Value = "Primera"
$ws.Cells.Item(2, 13).Value = 60
$ws.Cells.Item(2, 14).Value = 15000
$ws.Cells.Item(2, 15).Value = 15000
$ws.Cells.Item(2, 16).Value = 15000
$ws.Cells.Item(2, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(2, 18).Value = "Región Metropolitana"
$ws.Cells.Item(2, 19).Value = 833
$ws.Cells.Item(2, 20).Value = 18
